$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$r = $p.Range
$rStart = $d.Range($r.Start, $r.Start)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:spacing w:line="276" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t>Micrófono de bobina móvil:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve"> consta de un diafragma rígido suspendido frente a un imán permanente potente, que cuenta con una hendidura en la que va una bobina móvil acoplada</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:rPr>
            <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
            <w:bCs/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
            <w:lang w:val="es-CO"/>
          </w:rPr>
          <w:id w:val="411369749"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:instrText xml:space="preserve"> CITATION Vic15 \l 9226 </w:instrText>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:noProof/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:t>(Ruiz, 2015)</w:t>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve"> Las ondas de sonido provocan la vibración del diafragma (hacia adelante y hacia atrás) y a su vez de la bobina, el cambio de flujo magnético a través de las espiras de la bobina induce una corriente eléctrica en el conductor (conocido como ley de Faraday) proporcional a el sonido. Es de los más populares y baratos de construir. A nivel profesional no tiene la mejor fidelidad.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="276" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="276" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve">Micrófono de condensador: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t>consta de dos placas metálicas, una fija y otra móvil, separadas por un material aislante, similar a la estructura de un capacitor o condensador</w:t>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:rPr>
            <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
            <w:bCs/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
            <w:lang w:val="es-CO"/>
          </w:rPr>
          <w:id w:val="1958054860"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:instrText xml:space="preserve"> CITATION Vic15 \l 9226 </w:instrText>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:noProof/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:t xml:space="preserve"> </w:t>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:noProof/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:t>(Ruiz, 2015)</w:t>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve"> La placa que tiene libertad de movimiento vibrará con las ondas sonoras, y se alejará o acercará a la placa fija, esto provoca un cambio en la capacitancia del sistema (la capacitancia es inversamente proporcional a la distancia entre las placas). El cambio de la capacidad de almacenar carga es la responsable de producir la señal eléctrica. En este tipo de micrófono se necesita de un potencial o fuente de voltaje para funcionar (fuente fantasma). </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t>Generalmente c</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve">aros y frágiles. Gran fidelidad y calidad del sonido. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="276" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve">Micrófono piezoeléctrico: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t>consta de un diafragma y un elemento piezoeléctrico, compuesto de un material capaz de producir voltaje cuando se somete a presión mecánica</w:t>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:rPr>
            <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
            <w:bCs/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
            <w:lang w:val="es-CO"/>
          </w:rPr>
          <w:id w:val="-963421316"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:instrText xml:space="preserve"> CITATION Vic15 \l 9226 </w:instrText>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:noProof/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:t xml:space="preserve"> </w:t>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:noProof/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:t>(Ruiz, 2015)</w:t>
          </w:r>
          <w:r>
            <w:rPr>
              <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
              <w:bCs/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
              <w:lang w:val="es-CO"/>
            </w:rPr>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve"> Las ondas sonoras hacen vibrar el diafragma y, el movimiento de este hace que se mueva el material contenido en su interior (cuarzo, sal de Rochelle, carbón, etc.) </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>La fricción entre las partículas del material genera sobre la superficie de este una tensión eléctrica.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-CO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve">Hasta el momento solo se habló de algunos de los transductores o micrófonos que existen, sin embargo, el proceso de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t>sensar</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="es-CO"/>
        </w:rPr>
        <w:t xml:space="preserve"> sonido no se queda acá. Las señales eléctricas producidas por un micrófono suelen ser débiles, y requiere de amplificación, entre otros tratamientos de señales que apunte a las necesidades de la medición.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rStart.InsertXML($xml)
